$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the question text in A2 - reflow the line break
$ws.Range("A2").Value = "`"Governments should actively cooperate to have all countries`nconverge in terms of GDP per capita by the end of the century`""

# Reset row height after the text edit so no explicit custom height is persisted
$ws.Rows(2).AutoFit()

# Update numeric values in row 2 with final data
$ws.Range("B2").Value = 0.609601586795904
$ws.Range("K2").Value = 0.570384823995449
$ws.Range("L2").Value = 0.777916745185535
$ws.Range("N2").Value = 0.47407392499366
